$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.ClearFormats()

$ws.Range("A2").Value = "01‏/05‏/2025 02:07:49 م"
$ws.Range("B2").Value = "IDRF"
$ws.Range("C2").Value = "C2"
$ws.Range("D2").Value = "الرحلة 2"
$ws.Range("E2").Value = "الصمود"
$ws.Range("F2").Value = "أحمد شريم"

$ws.Range("G2").Value = "'1212"
$ws.Range("H2").Value = "'"

$null = $ws.Range("A1").Select()
